# Insert a new weekly record for "Femacal de La Calera" (Arándano/blue, Coquimbo)
# right above the existing row 107, shifting all subsequent rows down by one.
# This mirrors the source diff: dimension grows from A1:T157 to A1:T158 and a
# new row of data is inserted at row 107.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at position 107; Excel shifts rows 107:157 -> 108:158
# and extends the used range automatically.
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row 107 with the new record's values.
$ws.Range("A107").Value = 3
$ws.Range("B107").Value = "Femacal de La Calera"
$ws.Range("C107").Value = "Coquimbo"
$ws.Range("D107").Value = 44518
$ws.Range("E107").Value = 5
$ws.Range("F107").Value = "Fruta"
$ws.Range("G107").Value = 100101
$ws.Range("H107").Value = "Berries"
$ws.Range("I107").Value = 100101001
$ws.Range("J107").Value = "Arándano (blue)"
$ws.Range("K107").Value = "Sin especificar"
$ws.Range("L107").Value = "Primera"
$ws.Range("M107").Value = 58
$ws.Range("N107").Value = 15000
$ws.Range("O107").Value = 15000
$ws.Range("P107").Value = 15000
$ws.Range("Q107").Value = "$/bandeja 5 kilos"
$ws.Range("R107").Value = "Provincia de Curicó"
$ws.Range("S107").Value = 3000
$ws.Range("T107").Value = 5
